$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 125
$row2 = 126

$ws.Cells.Item($row1, 1).Value = 11
$ws.Cells.Item($row1, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row1, 3).Value = "Bíobío"
$ws.Cells.Item($row1, 4).Value = 44448
$ws.Cells.Item($row1, 4).NumberFormat = $ws.Cells.Item($row1 - 1, 4).NumberFormat
$ws.Cells.Item($row1, 5).Value = 8
$ws.Cells.Item($row1, 6).Value = 100114013
$ws.Cells.Item($row1, 7).Value = "Zanahoria"
$ws.Cells.Item($row1, 8).Value = "Sin especificar"
$ws.Cells.Item($row1, 9).Value = "Primera"
$ws.Cells.Item($row1, 10).Value = 600
$ws.Cells.Item($row1, 11).Value = 5000
$ws.Cells.Item($row1, 12).Value = 5500
$ws.Cells.Item($row1, 13).Value = 5250
$ws.Cells.Item($row1, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item($row1, 15).Value = "Región de Ñuble"
$ws.Cells.Item($row1, 16).Value = 262
$ws.Cells.Item($row1, 17).Value = 20
$ws.Cells.Item($row1, 18).Value = "Hortaliza"

$ws.Cells.Item($row2, 1).Value = 11
$ws.Cells.Item($row2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row2, 3).Value = "Bíobío"
$ws.Cells.Item($row2, 4).Value = 44448
$ws.Cells.Item($row2, 4).NumberFormat = $ws.Cells.Item($row1 - 1, 4).NumberFormat
$ws.Cells.Item($row2, 5).Value = 8
$ws.Cells.Item($row2, 6).Value = 100114013
$ws.Cells.Item($row2, 7).Value = "Zanahoria"
$ws.Cells.Item($row2, 8).Value = "Sin especificar"
$ws.Cells.Item($row2, 9).Value = "Segunda"
$ws.Cells.Item($row2, 10).Value = 300
$ws.Cells.Item($row2, 11).Value = 4500
$ws.Cells.Item($row2, 12).Value = 4500
$ws.Cells.Item($row2, 13).Value = 4500
$ws.Cells.Item($row2, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item($row2, 15).Value = "Región de Ñuble"
$ws.Cells.Item($row2, 16).Value = 225
$ws.Cells.Item($row2, 17).Value = 20
$ws.Cells.Item($row2, 18).Value = "Hortaliza"
